$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Convert the data range into an Excel Table (ListObject) named Table1 ---
# A1 is temporarily reset to the "Normal" style before the table is created so
# that Excel does not capture the header row's existing bold/fill/border
# formatting as a dedicated table "headerRowDxf" (the source file's header
# formatting is carried purely via the normal cell style, not a table dxf).
$ws.Range("A1").Style = "Normal"

$tableRange = $ws.Range("A1:U76")
$list = $ws.ListObjects.Add(1, $tableRange, $null, 1, $null)
$list.Name = "Table1"

# Restore A1's original header formatting by copying it from a neighboring
# header cell (whose formatting was untouched).
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 2. Rename header row (row 1) values: _old -> _FV2310, _new -> _FV2404 ---
# (The table's column names track these header cells automatically.)
for ($i = 1; $i -le 21; $i++) {
    $cell = $ws.Cells.Item(1, $i)
    $val = $cell.Value()
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2310"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2404"
        }
    }
}

# --- 3. Freeze the header row (first row) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
